$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.430167198181152
$ws.Range("B1").Value = 2.075509071350098
$ws.Range("C1").Value = 2.231804609298706
$ws.Range("D1").Value = 4.126387119293213
$ws.Range("E1").Value = 0.865297257900238
